$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.039.16'
$ws.Range("E2").Value = '  -0.89%  '
$ws.Range("D3").Value = '1.826.12'
$ws.Range("E3").Value = '  -0.33%  '
$ws.Range("E4").Value = '  -0.35%  '
$ws.Range("D5").Value = '''310.42'
$ws.Range("E5").Value = '  -1.32%  '
$ws.Range("E6").Value = '  -0.21%  '
$ws.Range("D7").Value = '''0.4623'
$ws.Range("E7").Value = '  -2.36%  '
$ws.Range("D8").Value = '''0.3692'
$ws.Range("E8").Value = '  +0.14%  '
$ws.Range("D9").Value = '''0.07252'
$ws.Range("E9").Value = '  -2.48%  '
$ws.Range("D10").Value = '''0.8618'
$ws.Range("E10").Value = '  -2.67%  '
$ws.Range("D11").Value = '''19.92'
$ws.Range("E11").Value = '  -2.75%  '
$ws.Range("D12").Value = '''0.07818'
$ws.Range("E12").Value = '  +6.72%  '
$ws.Range("E13").Value = '  -2.24%  '
$ws.Range("D14").Value = '''5.341'
$ws.Range("E14").Value = '  -1.61%  '
$ws.Range("D15").Value = '''6.536'
$ws.Range("E15").Value = '  -0.37%  '
$ws.Range("D16").Value = '''91.82'
$ws.Range("E16").Value = '  -2.36%  '
$ws.Range("E17").Value = '  -0.07%  '
$ws.Range("D18").Value = '''0.000008704'
$ws.Range("E18").Value = '  -0.96%  '
$ws.Range("E19").Value = '  -0.28%  '
$ws.Range("D20").Value = '27.161.52'
$ws.Range("D21").Value = '''14.52'
$ws.Range("E21").Value = '  -1.69%  '
$ws.Range("D22").Value = '''5.149'
$ws.Range("E22").Value = '  -2.53%  '
$ws.Range("E23").Value = '  -1.05%  '
$ws.Range("D24").Value = '2.084.28'
$ws.Range("E24").Value = '  -0.69%  '
$ws.Range("D25").Value = '''152.82'
$ws.Range("E25").Value = '  +0.64%  '
$ws.Range("E26").Value = '  -2.37%  '
$ws.Range("D27").Value = '''18.20'
$ws.Range("E27").Value = '  -2.31%  '
$ws.Range("D28").Value = '''2.093'
$ws.Range("E28").Value = '  -2.31%  '
$ws.Range("D29").Value = '''5.120'
$ws.Range("E29").Value = '  -2.03%  '
$ws.Range("D30").Value = '''115.43'
$ws.Range("E30").Value = '  -1.36%  '
$ws.Range("D31").Value = '''0.08840'
$ws.Range("E31").Value = '  -1.67%  '
$ws.Range("D32").Value = '''2.961'
$ws.Range("E32").Value = '  +0.71%  '
$ws.Range("D33").Value = '''4.438'
$ws.Range("E33").Value = '  -2.27%  '
$ws.Range("B34").Value = 'ImmutableX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D34").Value = '''0.7232'
$ws.Range("E34").Value = '  -3.53%  '
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").Value = '''1.133'
$ws.Range("E35").Value = '  -3.53%  '
$ws.Range("E36").Value = '  -0.63%  '
$ws.Range("D37").Value = '''2.450'
$ws.Range("E37").Value = '  +2.46%  '
$ws.Range("D38").Value = '''0.05240'
$ws.Range("E38").Value = '  -1.92%  '
$ws.Range("D39").Value = '''0.01941'
$ws.Range("E39").Value = '  -0.83%  '
$ws.Range("D40").Value = '''2.957'
$ws.Range("E40").Value = '  -0.75%  '
$ws.Range("D41").Value = '''7.227'
$ws.Range("E41").Value = '  -0.23%  '
$ws.Range("E42").Value = '  -2.45%  '
$ws.Range("E43").Value = '  -1.86%  '
$ws.Range("D44").Value = '''0.8593'
$ws.Range("E44").Value = '  -14.96%  '
$ws.Range("D45").Value = '''8.190'
$ws.Range("E45").Value = '  -3.43%  '
$ws.Range("D46").Value = '''0.4807'
$ws.Range("E46").Value = '  -2.49%  '
$ws.Range("D47").Value = '''1.008'
$ws.Range("E47").Value = '  -0.20%  '
$ws.Range("D48").Value = '''10.19'
$ws.Range("E48").Value = '  -3.42%  '
$ws.Range("D49").Value = '''102.74'
$ws.Range("E49").Value = '  -2.11%  '
$ws.Range("D50").Value = '''0.06265'
$ws.Range("E50").Value = '  -0.42%  '
$ws.Range("E51").Value = '  -2.94%  '
